# Generate Report for Handoff
#
# - Status text changes from "Handed back: in sync with en-US" to "Ready for handoff"
#   on the Overview, zh-cn and de-de sheets.
# - Timestamps are refreshed to reflect the new report generation time.
# - Column widths for the (now narrower) Status columns shrink to fit the
#   shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps
$overview.Range("G2").Value = "2016-09-04 09:03:30"
$dede.Range("H2").Value = "2016-09-04 09:03:30"
$zhcn.Range("H2").Value = "2016-09-04 09:03:26"

# --- Column widths (shrink Status columns now that text is shorter).
# Note: Excel's ColumnWidth setter snaps to the nearest displayable pixel
# width (quantized to 1/6 of a character for the default font), so the
# nearest reachable value to the recorded 17.2159881591797 is ~17.1667;
# 16.3 is the COM-visible width that lands on that pixel-quantized value.
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
